$wb = $excel.ActiveWorkbook

# --- Sheet "展览": update "想去人数" (column F) values ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 9118
$ws1.Cells.Item(3, 6).Value = 1984
$ws1.Cells.Item(4, 6).Value = 6648
$ws1.Cells.Item(6, 6).Value = 2168
$ws1.Cells.Item(7, 6).Value = 614
$ws1.Cells.Item(8, 6).Value = 86
$ws1.Cells.Item(9, 6).Value = 25
$ws1.Cells.Item(10, 6).Value = 83
$ws1.Cells.Item(13, 6).Value = 18
$ws1.Cells.Item(14, 6).Value = 92
$ws1.Cells.Item(15, 6).Value = 36
$ws1.Cells.Item(16, 6).Value = 9154
$ws1.Cells.Item(19, 6).Value = 209
$ws1.Cells.Item(21, 6).Value = 1865
$ws1.Cells.Item(23, 6).Value = 19
$ws1.Cells.Item(25, 6).Value = 115
$ws1.Cells.Item(28, 6).Value = 1054
$ws1.Cells.Item(29, 6).Value = 30
$ws1.Cells.Item(31, 6).Value = 576
$ws1.Cells.Item(33, 6).Value = 65
$ws1.Cells.Item(34, 6).Value = 560
$ws1.Cells.Item(35, 6).Value = 2458
$ws1.Cells.Item(36, 6).Value = 890
$ws1.Cells.Item(37, 6).Value = 569
$ws1.Cells.Item(40, 6).Value = 12
$ws1.Cells.Item(41, 6).Value = 332
$ws1.Cells.Item(43, 6).Value = 12
$ws1.Cells.Item(45, 6).Value = 38
$ws1.Cells.Item(46, 6).Value = 89
$ws1.Cells.Item(47, 6).Value = 33
$ws1.Cells.Item(48, 6).Value = 4011
$ws1.Cells.Item(49, 6).Value = 18

# --- Sheet "演出": remove 3 cancelled/duplicate rows (old rows 3-5), shift rows up ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows("3:5").Delete()
for ($i = 3; $i -le 26; $i++) {
    $ws2.Cells.Item($i, 1).Value = $i - 2
}

# --- Sheet "本地生活": update "想去人数" (column F) values ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 2349
$ws3.Cells.Item(4, 6).Value = 355

# --- Sheet "全部类型": update "想去人数" (column F) values ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 2349
$ws4.Cells.Item(3, 6).Value = 9118
$ws4.Cells.Item(5, 6).Value = 1984
$ws4.Cells.Item(6, 6).Value = 6648
$ws4.Cells.Item(8, 6).Value = 614
$ws4.Cells.Item(9, 6).Value = 86
$ws4.Cells.Item(10, 6).Value = 25
$ws4.Cells.Item(13, 6).Value = 83
$ws4.Cells.Item(15, 6).Value = 18
$ws4.Cells.Item(16, 6).Value = 92
$ws4.Cells.Item(17, 6).Value = 9155
$ws4.Cells.Item(20, 6).Value = 209
$ws4.Cells.Item(22, 6).Value = 1865
$ws4.Cells.Item(24, 6).Value = 115
$ws4.Cells.Item(27, 6).Value = 30
$ws4.Cells.Item(28, 6).Value = 4
$ws4.Cells.Item(29, 6).Value = 576
$ws4.Cells.Item(31, 6).Value = 65
$ws4.Cells.Item(32, 6).Value = 560
$ws4.Cells.Item(33, 6).Value = 890
$ws4.Cells.Item(36, 6).Value = 569
$ws4.Cells.Item(37, 6).Value = 332
$ws4.Cells.Item(42, 6).Value = 38
$ws4.Cells.Item(43, 6).Value = 89
$ws4.Cells.Item(44, 6).Value = 33
$ws4.Cells.Item(45, 6).Value = 4011
$ws4.Cells.Item(48, 6).Value = 18
